$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and B get a Text number format applied (matches the "Text" / numFmtId 49 style)
$ws.Range("A1:B5").NumberFormat = "@"

# New rows of user data
$ws.Range("A3").Value = "neena"
$ws.Range("B3").Value = "neena"

$ws.Range("A4").Value = "hana"
$ws.Range("B4").Value = "hana"

$ws.Range("A5").Value = '$yamala'
$ws.Range("B5").Value = "12@1234"

# Existing row 2's password value is updated
$ws.Range("B2").Value = "pppp"

# B5 becomes a hyperlink (auto-detected "@" style link), which also applies the Hyperlink style
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:12@1234")

# Selection moves to B2
$ws.Range("B2").Select() | Out-Null
